$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.764.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.56%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.827.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'598.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'166.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.53%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.828.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.13%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.96%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -4.67%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.45%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.40%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.20%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'36.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.44%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.466.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.830.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.54%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'67.909.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.24%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'18.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.85%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.56%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.93%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'464.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -7.11%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.729"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.59%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.0000160"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.65%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'82.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'2.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.68%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -3.40%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.41%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -3.52%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.13%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.973.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.17%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.27%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'30.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'9.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.21%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.785.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.30%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.75%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +9.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.13%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -3.13%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.04%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.311"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.12%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -7.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'FLOKI"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.000295"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.69%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Bittensor"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'418.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.47%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'8.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.97%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'47.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.75%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'141.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.01%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -3.79%  "
$ws.Range("E51").Style = "Normal"
